$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Figure out the sheet's current used range so we touch exactly the data
# columns (A..AH) and find the last two data rows (5 and 6).
$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count
$lastCol = $usedRange.Columns.Count

$targetRow = $lastRow - 1   # row 5: the row whose precision gets reduced
$dropRow   = $lastRow       # row 6: the row that gets removed entirely

# Reduce row 5 (B5:AH5) to 2 decimal places of precision - i.e. apply
# Excel's ROUND(value, 2) (round-half-away-from-zero) to every numeric
# reading cell in that row. Column A (the timestamp) is left untouched.
for ($col = 2; $col -le $lastCol; $col++) {
    $cell = $ws.Cells.Item($targetRow, $col)
    $orig = $cell.Value2
    if ($orig -ne $null) {
        $scaled = [double]$orig * 100.0
        if ($scaled -ge 0) {
            $rounded = [Math]::Floor($scaled + 0.5) / 100.0
        } else {
            $rounded = [Math]::Ceiling($scaled - 0.5) / 100.0
        }
        $cell.Value2 = $rounded
    }
}

# Remove the last row (row 6) entirely - the sheet keeps only one reading
# row after this edit, and the used range / dimension shrink accordingly.
$ws.Rows.Item($dropRow).Delete()
